$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the neighboring header cell (I2) onto the
# new header cell K2 before writing its value, so K2 picks up the same
# bold / wrap-text / vertical-center style already used by the other
# header cells (A2:I2) instead of the unstyled default.
$null = $ws.Range("I2").Copy()
$null = $ws.Range("K2").PasteSpecial(-4122)

# New column K: "IP Address 3" header plus five IP values.
$ws.Range("K2").Value = "IP Address 3"
$ws.Range("K3").Value = "13.232.89.90"
$ws.Range("K4").Value = "13.201.228.113"
$ws.Range("K5").Value = "3.109.3.143"
$ws.Range("K6").Value = "13.233.63.30"
$ws.Range("K7").Value = "3.109.59.156"

# Size the new column to fit its content.
$ws.Range("K1").ColumnWidth = 32.5

# The header row now wraps onto more lines, so it grows taller.
$ws.Range("A2:K2").RowHeight = 43.2

# Leave the selection where data entry would naturally finish.
$null = $ws.Range("K8").Select()
